$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.143.65"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "1.831.78"
$ws.Range("E3").Value = "  -3.25%  "
$ws.Range("D4").Value = "`'0.9997"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "`'230.34"
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("D6").Value = "`'0.9997"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "`'0.4640"
$ws.Range("E7").Value = "  -4.20%  "
$ws.Range("E8").Value = "  -6.36%  "
$ws.Range("D9").Value = "`'0.06251"
$ws.Range("E9").Value = "  -4.72%  "
$ws.Range("D10").Value = "1.826.92"
$ws.Range("E10").Value = "  -3.91%  "
$ws.Range("D11").Value = "`'0.07358"
$ws.Range("E11").Value = "  -1.42%  "
$ws.Range("D12").Value = "`'16.05"
$ws.Range("E12").Value = "  -4.92%  "
$ws.Range("D13").Value = "`'4.884"
$ws.Range("E13").Value = "  -4.50%  "
$ws.Range("D14").Value = "`'83.11"
$ws.Range("E14").Value = "  -5.71%  "
$ws.Range("D15").Value = "`'0.6192"
$ws.Range("E15").Value = "  -7.53%  "
$ws.Range("D16").Value = "30.084.51"
$ws.Range("E16").Value = "  -1.90%  "
$ws.Range("D17").Value = "`'0.9991"
$ws.Range("E17").Value = "  -0.08%  "
$ws.Range("D18").Value = "`'228.61"
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("D19").Value = "`'0.000007257"
$ws.Range("E19").Value = "  -4.28%  "
$ws.Range("D20").Value = "`'12.36"
$ws.Range("E20").Value = "  -6.84%  "
$ws.Range("D21").Value = "`'1.001"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("D22").Value = "2.066.93"
$ws.Range("E22").Value = "  -4.44%  "
$ws.Range("D23").Value = "`'4.831"
$ws.Range("E23").Value = "  -8.53%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "`'5.824"
$ws.Range("E24").Value = "  -6.26%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "`'165.20"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "`'9.079"
$ws.Range("E26").Value = "  -3.27%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "`'17.67"
$ws.Range("E27").Value = "  -6.46%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "`'1.840"
$ws.Range("E28").Value = "  -6.28%  "
$ws.Range("B29").Value = "Stellar"
$ws.Range("C29").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D29").Value = "`'0.1015"
$ws.Range("E29").Value = "  -0.94%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "`'1.366"
$ws.Range("E30").Value = "  -2.07%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "`'4.039"
$ws.Range("E31").Value = "  -7.04%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "`'3.748"
$ws.Range("E32").Value = "  -7.24%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "`'0.04781"
$ws.Range("E33").Value = "  -5.58%  "
$ws.Range("B34").Value = "ARBITRUM"
$ws.Range("C34").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D34").Value = "`'1.125"
$ws.Range("E34").Value = "  -7.43%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "`'0.6968"
$ws.Range("E35").Value = "  -7.63%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "`'2.686"
$ws.Range("E36").Value = "  -0.95%  "
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "`'0.01815"
$ws.Range("E37").Value = "  -3.52%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "`'2.608"
$ws.Range("E38").Value = "  -1.48%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "`'0.8925"
$ws.Range("E39").Value = "  -3.02%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "`'1.917"
$ws.Range("E40").Value = "  -7.43%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").Value = "`'0.9997"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").Value = "`'102.22"
$ws.Range("E42").Value = "  -4.68%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "`'5.483"
$ws.Range("E43").Value = "  -3.53%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").Value = "`'0.3991"
$ws.Range("E44").Value = "  -7.31%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "`'6.918"
$ws.Range("E45").Value = "  -6.96%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").Value = "`'0.1189"
$ws.Range("E46").Value = "  -6.97%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").Value = "`'59.35"
$ws.Range("E47").Value = "  -8.21%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "`'8.431"
$ws.Range("E48").Value = "  -5.90%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "`'0.05520"
$ws.Range("E49").Value = "  -2.59%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "`'32.45"
$ws.Range("E50").Value = "  -4.69%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "`'1.353"
$ws.Range("E51").Value = "  -9.55%  "
